$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "243.57" }
    @{ Cell = "G2"; Value = "14" }
    @{ Cell = "D3"; Value = "24.11" }
    @{ Cell = "G3"; Value = "14" }
    @{ Cell = "D4"; Value = "5.378" }
    @{ Cell = "G4"; Value = "14" }
    @{ Cell = "D5"; Value = "0.05920" }
    @{ Cell = "G5"; Value = "14" }
    @{ Cell = "D6"; Value = "3.393" }
    @{ Cell = "G6"; Value = "14" }
    @{ Cell = "D7"; Value = "6.500" }
    @{ Cell = "G7"; Value = "14" }
    @{ Cell = "D8"; Value = "0.8124" }
    @{ Cell = "G8"; Value = "14" }
    @{ Cell = "D9"; Value = "0.9512" }
    @{ Cell = "G9"; Value = "14" }
    @{ Cell = "G10"; Value = "14" }
    @{ Cell = "D11"; Value = "0.07414" }
    @{ Cell = "G11"; Value = "14" }
    @{ Cell = "D12"; Value = "0.03116" }
    @{ Cell = "G12"; Value = "14" }
    @{ Cell = "D13"; Value = "0.03063" }
    @{ Cell = "G13"; Value = "14" }
    @{ Cell = "D14"; Value = "0.09347" }
    @{ Cell = "G14"; Value = "14" }
    @{ Cell = "D15"; Value = "3.851" }
    @{ Cell = "G15"; Value = "14" }
    @{ Cell = "D16"; Value = "0.001576" }
    @{ Cell = "G16"; Value = "14" }
    @{ Cell = "D17"; Value = "0.04711" }
    @{ Cell = "G17"; Value = "14" }
    @{ Cell = "D18"; Value = "0.0005970" }
    @{ Cell = "E18"; Value = "17OneONE" }
    @{ Cell = "G18"; Value = "14" }
    @{ Cell = "D19"; Value = "0.005858" }
    @{ Cell = "G19"; Value = "14" }
    @{ Cell = "D20"; Value = "0.001247" }
    @{ Cell = "G20"; Value = "14" }
    @{ Cell = "G21"; Value = "14" }
    @{ Cell = "D22"; Value = "0.00009000" }
    @{ Cell = "G22"; Value = "14" }
    @{ Cell = "G23"; Value = "14" }
    @{ Cell = "G24"; Value = "14" }
    @{ Cell = "D25"; Value = "0.3224" }
    @{ Cell = "G25"; Value = "14" }
    @{ Cell = "G26"; Value = "14" }
    @{ Cell = "G27"; Value = "14" }
    @{ Cell = "G28"; Value = "14" }
    @{ Cell = "G29"; Value = "14" }
    @{ Cell = "G30"; Value = "14" }
    @{ Cell = "G31"; Value = "14" }
    @{ Cell = "G32"; Value = "14" }
    @{ Cell = "G33"; Value = "14" }
    @{ Cell = "G34"; Value = "14" }
    @{ Cell = "G35"; Value = "14" }
    @{ Cell = "G36"; Value = "14" }
    @{ Cell = "G37"; Value = "14" }
    @{ Cell = "G38"; Value = "14" }
    @{ Cell = "G39"; Value = "14" }
    @{ Cell = "D40"; Value = "0.03900" }
    @{ Cell = "G40"; Value = "14" }
    @{ Cell = "D41"; Value = "0.006340" }
    @{ Cell = "G41"; Value = "14" }
    @{ Cell = "G42"; Value = "14" }
    @{ Cell = "D43"; Value = "0.002840" }
    @{ Cell = "G43"; Value = "14" }
    @{ Cell = "D44"; Value = "0.008248" }
    @{ Cell = "G44"; Value = "14" }
    @{ Cell = "D45"; Value = "0.00005205" }
    @{ Cell = "G45"; Value = "14" }
    @{ Cell = "G46"; Value = "14" }
    @{ Cell = "D47"; Value = "0.6710" }
    @{ Cell = "G47"; Value = "14" }
    @{ Cell = "D48"; Value = "0.001996" }
    @{ Cell = "E48"; Value = "47BOLOBOLOWorstin24h" }
    @{ Cell = "G48"; Value = "14" }
    @{ Cell = "G49"; Value = "14" }
    @{ Cell = "G50"; Value = "14" }
    @{ Cell = "G51"; Value = "14" }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $col = $u.Cell.Substring(0, 1)
    if ($col -eq "D" -or $col -eq "G") {
        # These columns hold numeric-looking text in the source data;
        # force the Text number format so Excel keeps the literal
        # string instead of coercing it into a Number cell, then
        # restore the default style so formatting is unaffected.
        $range.NumberFormat = "@"
        $range.Value = $u.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $u.Value
    }
}
